$d = $word.ActiveDocument

$replacements = @(
    @("44×52=", "70×55="),
    @("61×27=", "78×50="),
    @("12×97=", "59×36="),
    @("47×58=", "22×18="),
    @("89×62=", "79×32="),
    @("38×89=", "73×87="),
    @("43×74=", "12×31="),
    @("93×65=", "25×14="),
    @("40×61=", "48×28="),
    @("36×14=", "46×77="),
    @("16×14=", "64×36="),
    @("20×78=", "63×49="),
    @("78×34=", "90×54="),
    @("37×17=", "88×55="),
    @("56×73=", "60×84="),
    @("41×48=", "14×90="),
    @("57×84=", "93×60="),
    @("35×86=", "27×78="),
    @("53×80=", "23×68="),
    @("14×35=", "19×66="),
    @("79×56=", "49×40="),
    @("17×62=", "33×73="),
    @("31×20=", "86×46="),
    @("64×45=", "97×41="),
    @("14×94=", "93×90=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
